# StorageComponentClassDiagram.pptx update:
#  - bump the cached "datetimeFigureOut" date placeholder text from 3/4/2019
#    to 4/14/2019 everywhere it appears (slide master, notes master, and
#    every slide layout)
#  - rename the "JsonAdaptedStatistics" class box to "JsonAdaptedRevenue"
#    on the (only) slide

$p = $ppt.ActivePresentation
$newDate = "4/14/2019"
$ppPlaceholderDate = 16

# ---------------------------------------------------------------------
# 1. Slide Master date placeholder
# ---------------------------------------------------------------------
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    $isDate = $false
    try {
        if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) { $isDate = $true }
    } catch {}
    if ($isDate) {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

# ---------------------------------------------------------------------
# 2. Notes Master date placeholder
# ---------------------------------------------------------------------
$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    $sh = $notesMaster.Shapes.Item($i)
    $isDate = $false
    try {
        if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) { $isDate = $true }
    } catch {}
    if ($isDate) {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

# ---------------------------------------------------------------------
# 3. Every Slide Layout's date placeholder
# ---------------------------------------------------------------------
$layouts = $master.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    $layout = $layouts.Item($L)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        $isDate = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) { $isDate = $true }
        } catch {}
        if ($isDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------------
# 4. Rename "JsonAdaptedStatistics" -> "JsonAdaptedRevenue" on slide 1
#    The shape lives inside a (possibly nested) group, so walk the group
#    tree up to two levels deep looking for the shape whose text is
#    "JsonAdaptedStatistics".
# ---------------------------------------------------------------------
$msoGroup = 6
$slide = $p.Slides.Item(1)
$target = $null

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $sh = $slide.Shapes.Item($i)

    if (($target -eq $null) -and $sh.HasTextFrame -and $sh.TextFrame.HasText) {
        if ($sh.TextFrame.TextRange.Text -eq "JsonAdaptedStatistics") {
            $target = $sh
        }
    }

    if (($target -eq $null) -and ($sh.Type -eq $msoGroup)) {
        $items1 = $sh.GroupItems
        for ($j = 1; $j -le $items1.Count; $j++) {
            $sh2 = $items1.Item($j)

            if (($target -eq $null) -and $sh2.HasTextFrame -and $sh2.TextFrame.HasText) {
                if ($sh2.TextFrame.TextRange.Text -eq "JsonAdaptedStatistics") {
                    $target = $sh2
                }
            }

            if (($target -eq $null) -and ($sh2.Type -eq $msoGroup)) {
                $items2 = $sh2.GroupItems
                for ($k = 1; $k -le $items2.Count; $k++) {
                    $sh3 = $items2.Item($k)
                    if ($sh3.HasTextFrame -and $sh3.TextFrame.HasText) {
                        if ($sh3.TextFrame.TextRange.Text -eq "JsonAdaptedStatistics") {
                            $target = $sh3
                        }
                    }
                }
            }
        }
    }
}

if ($target -ne $null) {
    $target.TextFrame.TextRange.Text = "JsonAdaptedRevenue"
}
